# Apply yellow highlighting to the "9. Hashing" section:
#  - "9. Hashing" heading: split into "9. " (unhighlighted) + "Hashing" (highlighted),
#    and the paragraph mark itself becomes highlighted.
#  - "Hash Tables: ..." bullet: both runs highlighted, paragraph mark highlighted too.
#  - "Applications of Hashing: ..." bullet: both runs highlighted (paragraph mark left as-is).

$d = $word.ActiveDocument

$wNs = 'xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml"'

# --- Paragraph: "9. Hashing" -------------------------------------------------
$range = $d.Content
$range.Find.Execute("9. Hashing") | Out-Null
$para = $range.Paragraphs(1)
$xml = '<w:p ' + $wNs + ' w14:paraId="08338B60" w14:textId="77777777" w:rsidR="00955902" w:rsidRPr="00955902" w:rsidRDefault="00955902" w:rsidP="00955902">' +
  '<w:pPr><w:rPr><w:b/><w:bCs/><w:highlight w:val="yellow"/></w:rPr></w:pPr>' +
  '<w:r w:rsidRPr="00955902"><w:rPr><w:b/><w:bCs/></w:rPr><w:t xml:space="preserve">9. </w:t></w:r>' +
  '<w:r w:rsidRPr="00955902"><w:rPr><w:b/><w:bCs/><w:highlight w:val="yellow"/></w:rPr><w:t>Hashing</w:t></w:r>' +
  '</w:p>'
$para.Range.InsertXML($xml) | Out-Null

# --- Paragraph: "Hash Tables: Hash functions, ..." ---------------------------
$range = $d.Content
$range.Find.Execute("Hash Tables:") | Out-Null
$para = $range.Paragraphs(1)
$xml = '<w:p ' + $wNs + ' w14:paraId="3C7693ED" w14:textId="77777777" w:rsidR="00955902" w:rsidRPr="00955902" w:rsidRDefault="00955902" w:rsidP="00955902">' +
  '<w:pPr><w:numPr><w:ilvl w:val="0"/><w:numId w:val="9"/></w:numPr><w:rPr><w:highlight w:val="yellow"/></w:rPr></w:pPr>' +
  '<w:r w:rsidRPr="00955902"><w:rPr><w:b/><w:bCs/><w:highlight w:val="yellow"/></w:rPr><w:t>Hash Tables:</w:t></w:r>' +
  '<w:r w:rsidRPr="00955902"><w:rPr><w:highlight w:val="yellow"/></w:rPr><w:t xml:space="preserve"> Hash functions, collision resolution techniques (chaining, open addressing)</w:t></w:r>' +
  '</w:p>'
$para.Range.InsertXML($xml) | Out-Null

# --- Paragraph: "Applications of Hashing: Frequency counting, ..." -----------
$range = $d.Content
$range.Find.Execute("Applications of Hashing:") | Out-Null
$para = $range.Paragraphs(1)
$xml = '<w:p ' + $wNs + ' w14:paraId="2C15D1D3" w14:textId="77777777" w:rsidR="00955902" w:rsidRPr="00955902" w:rsidRDefault="00955902" w:rsidP="00955902">' +
  '<w:pPr><w:numPr><w:ilvl w:val="0"/><w:numId w:val="9"/></w:numPr></w:pPr>' +
  '<w:r w:rsidRPr="00955902"><w:rPr><w:b/><w:bCs/><w:highlight w:val="yellow"/></w:rPr><w:t>Applications of Hashing:</w:t></w:r>' +
  '<w:r w:rsidRPr="00955902"><w:rPr><w:highlight w:val="yellow"/></w:rPr><w:t xml:space="preserve"> Frequency counting, finding duplicates, two-sum problem</w:t></w:r>' +
  '</w:p>'
$para.Range.InsertXML($xml) | Out-Null
